{"js": "const body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\n// \". Following allocution, Defendant\"  ->  \". Defendant\"\nconst change1 = body.search(\". Following allocution, Defendant\", { matchCase: true });\nchange1.load(\"text\");\nawait context.sync();\n\nif (change1.items.length > 0) {\n  change1.items[0].insertText(\". Defendant\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2 -------------------------------------------------------------\n// \"...accepted the plea and entered the following \"\n//   -> \"...accepted the plea and\" + \", following allocution,\" + \" entered the following \"\n// (moves \"following allocution\" from the first sentence into this clause,\n//  as three separate runs so the trailing run order/formatting matches.)\nconst change2 = body.search(\n  \"accepted the plea and entered the following \",\n  { matchCase: true }\n);\nchange2.load(\"text\");\nawait context.sync();\n\nif (change2.items.length > 0) {\n  change2.items[0].insertText(\n    \"accepted the plea and, following allocution, entered the following \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  // Force the newly-inserted \", following allocution,\" text to live in its\n  // own run (matching the target run layout) by nudging and then restoring\n  // a character property, which causes the host to split runs at the\n  // boundaries of the search hit instead of leaving one large merged run.\n  const mid = body.search(\", following allocution,\", { matchCase: true });\n  mid.load(\"text\");\n  await context.sync();\n\n  if (mid.items.length > 0) {\n    const midRange = mid.items[0];\n    midRange.font.bold = true;\n    await context.sync();\n    midRange.font.bold = false;\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 -------------------------------------------------------------\n# \". Following allocution, Defendant\"  ->  \". Defendant\"\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\". Following allocution, Defendant\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found1) {\n    $r1.Text = \". Defendant\"\n}\n\n# --- Change 2 ---------------------------------------------------------------\n# \"...accepted the plea and entered the following \"\n#   -> \"...accepted the plea and\" + \", following allocution,\" + \" entered the following \"\n# (moves \"following allocution\" out of the earlier sentence and into this\n#  clause, as three runs.)\n$r2 = $d.Content\n$found2 = $r2.Find.Execute(\"accepted the plea and entered the following \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found2) {\n    $r2.Text = \"accepted the plea and, following allocution, entered the following \"\n\n    # Re-find the newly inserted phrase and nudge a character property on it\n    # (set then clear Bold) so the host splits it into its own run instead of\n    # leaving it merged with the text before/after it.\n    $r3 = $d.Content\n    $found3 = $r3.Find.Execute(\", following allocution,\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    if ($found3) {\n        $r3.Font.Bold = 1\n        $r3.Font.Bold = 0\n    }\n}\n"}
